$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for columns D (label "U"), F and G for rows 2-11
$data = @(
    @(2, "U", 1, 1),
    @(3, "U", 3, 1),
    @(4, "U", 5, 1),
    @(5, "U", 3, 5),
    @(6, "U", 5, 5),
    @(7, "U", 6, 6),
    @(8, "U", 8, 6),
    @(9, "U", 10, 6),
    @(10, "U", 8, 10),
    @(11, "U", 10, 10)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value = $row[1]
    $ws.Cells.Item($r, 6).Value = $row[2]
    $ws.Cells.Item($r, 7).Value = $row[3]
}

# Update the selection to match the edited range
$null = $ws.Range("D2:G11").Select()
